$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value()
    if ($val -match '^(\d{4}-\d{2})-07$') {
        $newVal = "$($matches[1])-09"
        $cell.NumberFormat = "@"
        $cell.Value = $newVal
        $cell.Style = "Normal"
    }
}
